$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("Summary")
$wsZonal = $wb.Worksheets.Item("Zonal Data")
$wsTrades = $wb.Worksheets.Item("Trades")

# Summary sheet updates
$wsSummary.Range("B2").Value = 400
$wsSummary.Range("C2").Value = -400

# Zonal Data sheet updates
$wsZonal.Range("E2").Value = "X"
$wsZonal.Range("D3").Value = -150
$wsZonal.Range("E3").Value = "X"
$wsZonal.Range("C4").Value = 150
$wsZonal.Range("E4").Value = "X"

# Trades sheet updates
$wsTrades.Range("D2").Value = 150
$wsTrades.Range("D3").Value = 50
